$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column E to match new data (COM ColumnWidth has a +5/6 char offset
# vs. the stored OOXML <col width>, so back that out to land on exactly 25)
$ws.Columns.Item(5).ColumnWidth = 24.1666666666667

# Add the new row of data (row 4)
$ws.Range("A4").Value = "f5 address 25"
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = "f5 first 25"
$ws.Range("D4").Value = "f5 last 25"
$ws.Range("E4").Value = "{{credit_debit_number}}"
$ws.Range("F4").Value = "f5 city 25"

# G4/H4 look numeric ("2505" / "250505") but must stay text, matching the
# rest of the sheet (e.g. G2/G3). Force text storage, then drop the
# NumberFormat-derived style so no new cell style is introduced.
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2505"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "250505"
$ws.Range("G4:H4").Style = "Normal"
